# Progress Update 3 PPT
# - Refresh the cached "datetimeFigureOut" fields on the Handout Master and
#   Notes Master (11/23/2021 -> 11/27/2021).
# - Add "Thank You" to the title placeholder on the last (closing) slide.
# - Nudge the pie-chart group on the "Next Steps" slide down slightly.

$p = $ppt.ActivePresentation

# --- 1) Handout Master date placeholder -------------------------------
$handoutMaster = $p.HandoutMaster
$handoutMaster.HeadersFooters.DateAndTime.Text = "11/27/2021"

# --- 2) Notes Master date placeholder ----------------------------------
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "11/27/2021"

# --- 3) Closing slide: add "Thank You" to the title placeholder --------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$titleShape = $lastSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Thank You"
$titleShape.TextFrame.TextRange.LanguageID = "en-IN"

# --- 4) "Next Steps" slide: reposition the pie-chart group -------------
# (Target offset is 2251275 EMU; expressed directly in points, with a few
# extra significant digits, so the float round-trip back to EMU lands on
# the exact target instead of being truncated a hair short of it.)
$stepsSlide = $p.Slides.Item(9)
$group = $stepsSlide.Shapes.Item(5)
$group.Top = 177.26575503149607
